$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.146.27"
$ws.Range("E2").Value = "  -1.10%  "
$ws.Range("D3").Value = "3.065.01"
$ws.Range("E3").Value = "  -1.66%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'572.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.84%  "
$ws.Range("D6").Value = "'169.36"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.24%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "3.062.67"
$ws.Range("E8").Value = "  -1.54%  "
$ws.Range("D9").Value = "'0.509"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.25%  "
$ws.Range("D10").Value = "'6.28"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.12%  "
$ws.Range("E11").Value = "  -3.05%  "
$ws.Range("E12").Value = "  -3.24%  "
$ws.Range("E13").Value = "  -4.28%  "
$ws.Range("D14").Value = "'35.55"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.37%  "
$ws.Range("E15").Value = "  -1.64%  "
$ws.Range("D16").Value = "3.574.92"
$ws.Range("E16").Value = "  -1.59%  "
$ws.Range("D17").Value = "66.046.86"
$ws.Range("E17").Value = "  -1.20%  "
$ws.Range("D18").Value = "'6.93"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.91%  "
$ws.Range("D19").Value = "3.061.91"
$ws.Range("E19").Value = "  -1.65%  "
$ws.Range("D20").Value = "'16.43"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.40%  "
$ws.Range("D21").Value = "'483.92"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.86%  "
$ws.Range("D22").Value = "'0.684"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.58%  "
$ws.Range("D23").Value = "'7.62"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.68%  "
$ws.Range("D24").Value = "'82.24"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.91%  "
$ws.Range("D25").Value = "'12.60"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.80%  "
$ws.Range("E26").Value = "  -3.73%  "
$ws.Range("D27").Value = "'10.14"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.15%  "
$ws.Range("D28").Value = "'1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.01%  "
$ws.Range("D29").Value = "'7.85"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.10%  "
$ws.Range("D30").Value = "'2.24"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.47%  "
$ws.Range("D31").Value = "'2.59"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.33%  "
$ws.Range("D32").Value = "'27.62"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.14%  "
$ws.Range("E33").Value = "  -3.92%  "
$ws.Range("D34").Value = "0.0₃0898"
$ws.Range("E34").Value = "  -5.59%  "
$ws.Range("D35").Value = "'0.999"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.11%  "
$ws.Range("D36").Value = "'47.29"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.70%  "
$ws.Range("D37").Value = "'0.939"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.66%  "
$ws.Range("D38").Value = "'5.53"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.23%  "
$ws.Range("E39").Value = "  -1.41%  "
$ws.Range("D40").Value = "'1.95"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.71%  "
$ws.Range("D41").Value = "'0.298"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.25%  "
$ws.Range("D42").Value = "'8.20"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -5.30%  "
$ws.Range("D43").Value = "2.759.87"
$ws.Range("E43").Value = "  -1.84%  "
$ws.Range("D44").Value = "'2.51"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.21%  "
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").Value = "'0.0342"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.39%  "
$ws.Range("B46").Value = "Monero"
$ws.Range("C46").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D46").Value = "'134.28"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.62%  "
$ws.Range("D47").Value = "'361.74"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.46%  "
$ws.Range("D49").Value = "'24.16"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.72%  "
$ws.Range("E50").Value = "  -2.34%  "
$ws.Range("E51").Value = "  -3.01%  "
